# Apply marksheet corrections: update correct/total marks
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Row 11 "Marking" - B11 changes from 3 to 5
$ws.Range("B11").Value = 5

# Row 12 "Total" - B12 changes from 69 to 115, and E12 text changes from "66/84" to "115/140"
$ws.Range("B12").Value = 115
$ws.Range("E12").Value = "115/140"
